$d = $word.ActiveDocument

# The document's single table holds the HbA1c control-status rows in
# column 1 ("Controlled (HbA1c<7.5%)" / "Uncontrolled (HbA1c>7.5%)").
# Update the unit text from "7.5%" to "58mmol/mol" (and the ">"
# comparator to "≥") while leaving everything else (the leading
# "(HbA1c<" / "(HbA1c" text and the trailing ")") untouched. The
# shipped XML splits the cell text into three runs with identical run
# properties (prefix / middle "unit" chunk / trailing ")"), so each
# edit below is done in two steps: first a plain text replace that
# keeps the paragraph as a single run, and then a harmless format
# toggle on the inserted chunk that forces it to stay a distinct run
# (matching the shape of the shipped diff) instead of being re-merged
# with its neighbours.

function Find-RowWithPrefix {
    param(
        [string]$Prefix
    )

    $table = $d.Tables.Item(1)
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        $cellText = $table.Cell($r, 1).Range.Text
        if ($cellText.StartsWith($Prefix)) {
            return $r
        }
    }
    throw "Find-RowWithPrefix: no row found whose first cell starts with '$Prefix'"
}

function Update-HbA1cUnit {
    param(
        [string]$RowPrefix,
        [string]$OldUnitText,
        [string]$NewUnitText
    )

    $rowIndex = Find-RowWithPrefix $RowPrefix

    $table = $d.Tables.Item(1)
    $cell = $table.Cell($rowIndex, 1)
    $cellRange = $cell.Range
    $cellStart = $cellRange.Start
    $cellText = $cellRange.Text

    $unitStart = $cellText.IndexOf($OldUnitText)
    if ($unitStart -lt 0) {
        throw "Update-HbA1cUnit: could not find '$OldUnitText' in row $rowIndex (text was '$cellText')"
    }

    $rangeStart = $cellStart + $unitStart
    $rangeEnd = $rangeStart + $OldUnitText.Length

    $unitRange = $d.Range($rangeStart, $rangeEnd)
    $unitRange.Text = $NewUnitText

    # Re-fetch the now-stale range at the freshly inserted text's
    # position and nudge its formatting off/on so it is preserved as
    # its own run rather than being coalesced back into the runs
    # before/after it.
    $newRangeEnd = $rangeStart + $NewUnitText.Length
    $insertedRange = $d.Range($rangeStart, $newRangeEnd)
    $insertedRange.Font.Bold = $true
    $insertedRange.Font.Bold = $false
}

Update-HbA1cUnit "Controlled (HbA1c<" "7.5%" "58mmol/mol"
Update-HbA1cUnit "Uncontrolled (HbA1c" ">7.5%" "≥58mmol/mol"
